$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 360
$ws.Range("I6").Value = 332
$ws.Range("K6").Value = 996
$ws.Range("M6").Value = -884
$ws.Range("H19").Value = 585.64
$ws.Range("I19").Value = 575.3
$ws.Range("K19").Value = 575.3
$ws.Range("M19").Value = -400.3
$ws.Range("H29").Value = 3383.8333
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 3383.8333
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 10151.4999
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -10713.4999
$ws.Range("H32").Value = 4646663
$ws.Range("I32").Value = 219.5
$ws.Range("J32").Value = 6336279
$ws.Range("K32").Value = 219.5
$ws.Range("L32").Value = 6336279
$ws.Range("M32").Value = 106.5
$ws.Range("N32").Value = -6336931
$ws.Range("H33").Value = 697.4358999999999
$ws.Range("I33").Value = 438.26666
$ws.Range("J33").Value = 1561.3334
$ws.Range("K33").Value = 438.26666
$ws.Range("L33").Value = 1561.3334
$ws.Range("M33").Value = -209.26666
$ws.Range("N33").Value = -2019.3334
$ws.Range("H58").Value = 2161.35
$ws.Range("J58").Value = 2892.6428
$ws.Range("L58").Value = 8677.928400000001
$ws.Range("N58").Value = -8977.928400000001
$ws.Range("H86").Value = 3585
$ws.Range("I86").Value = 1983.2858
$ws.Range("J86").Value = 5563.5884
$ws.Range("K86").Value = 1983.2858
$ws.Range("L86").Value = 5563.5884
$ws.Range("M86").Value = -860.2858000000001
$ws.Range("N86").Value = -7809.5884
$ws.Range("H89").Value = 3585
$ws.Range("I89").Value = 1983.2858
$ws.Range("J89").Value = 5563.5884
$ws.Range("K89").Value = 9916.429
$ws.Range("L89").Value = 27817.942
$ws.Range("M89").Value = -4300.429
$ws.Range("N89").Value = -39049.942
$ws.Range("H116").Value = 120588.336
$ws.Range("I116").Value = 165199.23
$ws.Range("J116").Value = 4600
$ws.Range("K116").Value = 165199.23
$ws.Range("L116").Value = 4600
$ws.Range("M116").Value = -161757.23
$ws.Range("N116").Value = -11484
$ws.Range("H132").Value = 2745.8904
$ws.Range("I132").Value = 1470.6333
$ws.Range("J132").Value = 8631.691999999999
$ws.Range("K132").Value = 4411.8999
$ws.Range("L132").Value = 25895.076
$ws.Range("M132").Value = -1881.8999
$ws.Range("N132").Value = -30955.076
$ws.Range("H141").Value = 1786.3103
$ws.Range("I141").Value = 1763.2222
$ws.Range("J141").Value = 2098
$ws.Range("K141").Value = 5289.6666
$ws.Range("L141").Value = 6294
$ws.Range("M141").Value = -109.6665999999996
$ws.Range("N141").Value = -16654

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1128.0857
$ws.Range("I45").Value = 1004.13635
$ws.Range("J45").Value = 1337.8462
$ws.Range("K45").Value = 1004.13635
$ws.Range("L45").Value = 1337.8462
$ws.Range("M45").Value = -627.13635
$ws.Range("N45").Value = -2091.8462
$ws.Range("H61").Value = 3077.7046
$ws.Range("I61").Value = 1911.5807
$ws.Range("J61").Value = 5858.4614
$ws.Range("K61").Value = 1911.5807
$ws.Range("L61").Value = 5858.4614
$ws.Range("M61").Value = -1699.5807
$ws.Range("N61").Value = -6282.4614
$ws.Range("H110").Value = 1323.7858
$ws.Range("I110").Value = 1177.75
$ws.Range("J110").Value = 2200
$ws.Range("K110").Value = 1177.75
$ws.Range("L110").Value = 2200
$ws.Range("M110").Value = 867.25
$ws.Range("N110").Value = -6290
$ws.Range("H136").Value = 3077.7046
$ws.Range("I136").Value = 1911.5807
$ws.Range("J136").Value = 5858.4614
$ws.Range("K136").Value = 5734.742099999999
$ws.Range("L136").Value = 17575.3842
$ws.Range("M136").Value = -3184.742099999999
$ws.Range("N136").Value = -22675.3842

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3462.5293
$ws.Range("I31").Value = 2486.1667
$ws.Range("J31").Value = 4330.407
$ws.Range("K31").Value = 2486.1667
$ws.Range("L31").Value = 4330.407
$ws.Range("M31").Value = -2191.1667
$ws.Range("N31").Value = -4920.407
$ws.Range("H34").Value = 3462.5293
$ws.Range("I34").Value = 2486.1667
$ws.Range("J34").Value = 4330.407
$ws.Range("K34").Value = 2486.1667
$ws.Range("L34").Value = 4330.407
$ws.Range("M34").Value = -2284.1667
$ws.Range("N34").Value = -4734.407
$ws.Range("H58").Value = 2206.0789
$ws.Range("I58").Value = 1625.4736
$ws.Range("J58").Value = 2786.6843
$ws.Range("K58").Value = 1625.4736
$ws.Range("L58").Value = 2786.6843
$ws.Range("M58").Value = -1422.4736
$ws.Range("N58").Value = -3192.6843
$ws.Range("H105").Value = 610.6094000000001
$ws.Range("J105").Value = 661.25
$ws.Range("L105").Value = 661.25
$ws.Range("N105").Value = -4155.25
$ws.Range("H107").Value = 372.72223
$ws.Range("I107").Value = 284.23077
$ws.Range("J107").Value = 602.8
$ws.Range("K107").Value = 284.23077
$ws.Range("L107").Value = 602.8
$ws.Range("M107").Value = 1635.76923
$ws.Range("N107").Value = -4442.8
$ws.Range("H132").Value = 2308.52
$ws.Range("I132").Value = 1400.4517
$ws.Range("J132").Value = 3790.1052
$ws.Range("K132").Value = 4201.355100000001
$ws.Range("L132").Value = 11370.3156
$ws.Range("M132").Value = -1671.355100000001
$ws.Range("N132").Value = -16430.3156
$ws.Range("H136").Value = 2206.0789
$ws.Range("I136").Value = 1625.4736
$ws.Range("J136").Value = 2786.6843
$ws.Range("K136").Value = 4876.4208
$ws.Range("L136").Value = 8360.052899999999
$ws.Range("M136").Value = -2326.4208
$ws.Range("N136").Value = -13460.0529

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 102418.6
$ws.Range("I22").Value = 1477.5
$ws.Range("J22").Value = 169712.67
$ws.Range("K22").Value = 4432.5
$ws.Range("L22").Value = 509138.01
$ws.Range("M22").Value = -4263.5
$ws.Range("N22").Value = -509476.01
$ws.Range("H27").Value = 102418.6
$ws.Range("I27").Value = 1477.5
$ws.Range("J27").Value = 169712.67
$ws.Range("K27").Value = 4432.5
$ws.Range("L27").Value = 509138.01
$ws.Range("M27").Value = -4330.5
$ws.Range("N27").Value = -509342.01
$ws.Range("H131").Value = 1411.2424
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 1411.2424
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 4233.7272
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -14313.7272
$ws.Range("H132").Value = 5275
$ws.Range("I132").Value = 1485.7142
$ws.Range("J132").Value = 8222.223
$ws.Range("K132").Value = 13371.4278
$ws.Range("L132").Value = 74000.007
$ws.Range("M132").Value = -10841.4278
$ws.Range("N132").Value = -79060.007

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2575
$ws.Range("I102").Value = 2718.9395
$ws.Range("J102").Value = 1783.3334
$ws.Range("K102").Value = 2718.9395
$ws.Range("L102").Value = 1783.3334
$ws.Range("M102").Value = -1096.9395
$ws.Range("N102").Value = -5027.3334
$ws.Range("H132").Value = 3000.982
$ws.Range("I132").Value = 2751.3022
$ws.Range("J132").Value = 3895.6667
$ws.Range("K132").Value = 8253.9066
$ws.Range("L132").Value = 11687.0001
$ws.Range("M132").Value = -5723.9066
$ws.Range("N132").Value = -16747.0001
$ws.Range("H133").Value = 32853.332
$ws.Range("J133").Value = 32853.332
$ws.Range("L133").Value = 32853.332
$ws.Range("N133").Value = -42973.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 4400
$ws.Range("I5").Value = 3995
$ws.Range("J5").Value = 4670
$ws.Range("K5").Value = 3995
$ws.Range("L5").Value = 4670
$ws.Range("M5").Value = -3882
$ws.Range("N5").Value = -4896
$ws.Range("H132").Value = 5270.5737
$ws.Range("I132").Value = 1783.5385
$ws.Range("J132").Value = 11452.137
$ws.Range("K132").Value = 5350.6155
$ws.Range("L132").Value = 34356.411
$ws.Range("M132").Value = -2820.6155
$ws.Range("N132").Value = -39416.411
$ws.Range("H136").Value = 5161.3335
$ws.Range("I136").Value = 2768.7273
$ws.Range("K136").Value = 8306.1819
$ws.Range("M136").Value = -5756.1819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 206
$ws.Range("I113").Value = 214.83333
$ws.Range("J113").Value = 100
$ws.Range("K113").Value = 644.49999
$ws.Range("L113").Value = 300
$ws.Range("M113").Value = 1525.50001
$ws.Range("N113").Value = -4640
$ws.Range("H132").Value = 2057.7415
$ws.Range("I132").Value = 1153.4595
$ws.Range("J132").Value = 3651
$ws.Range("K132").Value = 3460.3785
$ws.Range("L132").Value = 10953
$ws.Range("M132").Value = -930.3784999999998
$ws.Range("N132").Value = -16013
$ws.Range("H136").Value = 1625.5068
$ws.Range("I136").Value = 893.7
$ws.Range("J136").Value = 2512.5454
$ws.Range("K136").Value = 2681.1
$ws.Range("L136").Value = 7537.6362
$ws.Range("M136").Value = -131.1000000000004
$ws.Range("N136").Value = -12637.6362
